$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 78, shifting existing rows 78-99 down to 79-100.
$ws.Rows.Item(78).Insert()

# Populate the newly inserted row 78 with the new weekly record.
$ws.Range("A78").Value = 7
$ws.Range("B78").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C78").Value = 'Ñuble'
$ws.Range("D78").Value = 44841
$ws.Range("E78").Value = 16
$ws.Range("F78").Value = 100112031
$ws.Range("G78").Value = 'Poroto verde'
$ws.Range("H78").Value = 'Magnum'
$ws.Range("I78").Value = 'Primera'
$ws.Range("J78").Value = 60
$ws.Range("K78").Value = 26000
$ws.Range("L78").Value = 27000
$ws.Range("M78").Value = 26500
$ws.Range("N78").Value = '$/malla 25 kilos'
$ws.Range("O78").Value = 'Perú'
$ws.Range("P78").Value = 1060
$ws.Range("Q78").Value = 25
$ws.Range("R78").Value = 'Hortaliza'
